$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.921134471893311
$ws.Range("B1").Value = 2.543084621429443
$ws.Range("C1").Value = 1.869756460189819
$ws.Range("D1").Value = 1.735360026359558
$ws.Range("E1").Value = 1.625322103500366
